$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "rwEac954"
$ws.Range("B2").Value = 23110323
$ws.Range("C2").Value = "mcjpkik97"
$ws.Range("D2").Value = "zyH`$#V85"
$ws.Range("F2").Value = "gOQrfKSa"
$ws.Range("G2").Value = "HKNd"

# Row 3
$ws.Range("A3").Value = "yduIu918"
$ws.Range("B3").Value = 23110322
$ws.Range("C3").Value = "zsazmox88"
$ws.Range("D3").Value = "Vs2T#6b`$"
$ws.Range("F3").Value = "HafZXSaF"
$ws.Range("G3").Value = "tOKJ"

# Row 4
$ws.Range("A4").Value = "PFUFu512"
$ws.Range("B4").Value = 23110321
$ws.Range("C4").Value = "dbxhczi18"
$ws.Range("D4").Value = "S&kK%e97"
$ws.Range("F4").Value = "qhBmTjfc"
$ws.Range("G4").Value = "MVSF"
